$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.969.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.014.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.56%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.564"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.026.49"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.23%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("E11").Value = "  -4.14%  "
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.554.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.031.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.024.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "395.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("E22").Value = "  -3.89%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.466"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("E26").Value = "  -4.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0976"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.15%  "
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "159.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.516.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.668"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0598"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.93%  "
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.52%  "
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "264.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.14%  "
